$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Classic_Chocolate / rahul ---
$ws.Range("D2").Value = 2
$ws.Range("F2").Value = 60
$ws.Range("G2").Value = 45283.51190972222

# --- Row 3: Vanilla / rahul ---
$ws.Range("D3").Value = 12
$ws.Range("F3").Value = 480
$ws.Range("G3").Value = 45283.51202546297

# --- Row 4: Bliss, order_id 2->1, customer raj->rahul ---
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "rahul"
$ws.Range("D4").Value = 10
$ws.Range("F4").Value = 350
$ws.Range("G4").Value = 45283.51260416667

# --- Row 5: order_id 2->3, customer raj->rohit, cake Cookies->Red_Velvet ---
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "rohit"
$ws.Range("C5").Value = "Red_Velvet"
$ws.Range("D5").Value = 10
$ws.Range("E5").Value = 45
$ws.Range("F5").Value = 450
$ws.Range("G5").Value = 45283.51989583333

# --- Row 6: order_id 1->4, customer rahul->raj, cake Cookies->Classic_Chocolate ---
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "raj"
$ws.Range("C6").Value = "Classic_Chocolate"
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = 30
$ws.Range("F6").Value = 150
$ws.Range("G6").Value = 45283.52236111111

# --- New row 7: raj / Bliss ---
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "raj"
$ws.Range("C7").Value = "Bliss"
$ws.Range("D7").Value = 8
$ws.Range("E7").Value = 35
$ws.Range("F7").Value = 280
$ws.Range("G7").Value = 45283.52288194445
$ws.Range("G7").NumberFormat = $ws.Range("G2").NumberFormat

# --- New row 8: raj / Vanilla ---
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "raj"
$ws.Range("C8").Value = "Vanilla"
$ws.Range("D8").Value = 22
$ws.Range("E8").Value = 40
$ws.Range("F8").Value = 880
$ws.Range("G8").Value = 45283.52310185185
$ws.Range("G8").NumberFormat = $ws.Range("G2").NumberFormat

# --- New row 9: raj / Vanilla ---
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "raj"
$ws.Range("C9").Value = "Vanilla"
$ws.Range("D9").Value = 55
$ws.Range("E9").Value = 40
$ws.Range("F9").Value = 2200
$ws.Range("G9").Value = 45283.52400462963
$ws.Range("G9").NumberFormat = $ws.Range("G2").NumberFormat

Write-Host "edits applied"
